$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: rank stays the same (column A), only the live-quote
# columns (Coin/Link/Price/Volume(1h)) move. Numeric-looking Price strings
# get a leading apostrophe so Excel keeps them as text (matches the sheet's
# existing "nnn.nnn.nn"-style text storage) instead of coercing to a number.
$updates = @(
    @{ "D2" = '60.593.90'; "E2" = '  -3.31%  ' }
    @{ "D3" = '3.346.57'; "E3" = '  -2.86%  ' }
    @{ "D4" = '''0.999'; "E4" = '  +0.00%  ' }
    @{ "D5" = '''566.81'; "E5" = '  -2.32%  ' }
    @{ "D6" = '''147.16'; "E6" = '  -0.65%  ' }
    @{ "E7" = '  -0.08%  ' }
    @{ "D8" = '''0.483'; "E8" = '  +0.24%  ' }
    @{ "D9" = '''7.93'; "E9" = '  -0.74%  ' }
    @{ "E10" = '  -1.22%  ' }
    @{ "D11" = '''0.415'; "E11" = '  +1.38%  ' }
    @{ "D12" = '3.913.80'; "E12" = '  -3.05%  ' }
    @{ "E13" = '  +1.18%  ' }
    @{ "D14" = '''27.74'; "E14" = '  -2.04%  ' }
    @{ "D15" = '3.340.04'; "E15" = '  -3.47%  ' }
    @{ "E16" = '  -1.54%  ' }
    @{ "D17" = '60.618.37'; "E17" = '  -3.34%  ' }
    @{ "D18" = '''6.27'; "E18" = '  -1.25%  ' }
    @{ "D19" = '''14.56'; "E19" = '  -0.51%  ' }
    @{ "D20" = '''8.91'; "E20" = '  -1.62%  ' }
    @{ "D21" = '''376.25'; "E21" = '  -2.63%  ' }
    @{ "D22" = '''0.560'; "E22" = '  -0.34%  ' }
    @{ "D23" = '''74.75'; "E23" = '  -0.75%  ' }
    @{ "E24" = '  +0.10%  ' }
    @{ "D25" = '3.491.12'; "E25" = '  -2.52%  ' }
    @{ "E26" = '  -5.50%  ' }
    @{ "E27" = '  -4.41%  ' }
    @{ "E28" = '  +0.14%  ' }
    @{ "D29" = '''7.35'; "E29" = '  -4.16%  ' }
    @{ "E30" = '  -1.38%  ' }
    @{ "E31" = '  +0.02%  ' }
    @{ "D32" = '''7.72'; "E32" = '  -3.59%  ' }
    @{ "E33" = '  -1.54%  ' }
    @{ "D34" = '''1.30'; "E34" = '  -2.89%  ' }
    @{ "D35" = '''5.32'; "E35" = '  -0.78%  ' }
    @{ "E36" = '  -4.86%  ' }
    @{ "D37" = '''6.83'; "E37" = '  -1.74%  ' }
    @{ "D38" = '''167.52'; "E38" = '  -1.15%  ' }
    @{ "E39" = '  -12.92%  ' }
    @{ "D40" = '3.379.27'; "E40" = '  -2.99%  ' }
    @{ "D41" = '''0.0750'; "E41" = '  -3.17%  ' }
    @{ "E42" = '  -3.64%  ' }
    @{ "E43" = '  -1.52%  ' }
    @{ "E44" = '  -2.94%  ' }
    @{ "E45" = '  -4.92%  ' }
    @{ "D46" = '2.456.40'; "E46" = '  -4.39%  ' }
    @{ "E47" = '  -3.29%  ' }
    @{ "E48" = '  -0.04%  ' }
    @{ "D49" = '''22.41'; "E49" = '  -1.10%  ' }
    @{ "E50" = '  -1.92%  ' }
    @{ "B51" = 'dogwifhat'; "C51" = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; "D51" = '''2.06'; "E51" = '  -8.04%  ' }
)

foreach ($update in $updates) {
    foreach ($cellRef in $update.Keys) {
        $ws.Range($cellRef).Value = $update[$cellRef]
    }
}
